$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for row 9 (previously blank D9:F9) ---
$ws.Range("D9").Value = 182819906
$ws.Range("E9").Value = 177151873
$ws.Range("F9").Value = 177683557

# --- Conditional formatting: 2-colour scale + 3-colour scale on row 9, then 3-colour
#     scales on rows 8 down to 3 (matches the document order / priority numbering
#     produced by the original edit) ---
$cf9b = $ws.Range("A9:XFD9").FormatConditions.AddColorScale(2)
$cf9  = $ws.Range("A9:XFD9").FormatConditions.AddColorScale(3)
$cf8  = $ws.Range("A8:XFD8").FormatConditions.AddColorScale(3)
$cf7  = $ws.Range("A7:XFD7").FormatConditions.AddColorScale(3)
$cf6  = $ws.Range("A6:XFD6").FormatConditions.AddColorScale(3)
$cf5  = $ws.Range("A5:XFD5").FormatConditions.AddColorScale(3)
$cf4  = $ws.Range("A4:XFD4").FormatConditions.AddColorScale(3)
$cf3  = $ws.Range("A3:XFD3").FormatConditions.AddColorScale(3)

# Re-assign priorities to match the target numbering (does not change file order,
# only the priority attribute each rule is saved with).
$cf9b.Priority = 8
$cf9.Priority  = 7
$cf8.Priority  = 6
$cf7.Priority  = 5
$cf6.Priority  = 4
$cf5.Priority  = 3
$cf4.Priority  = 2
$cf3.Priority  = 1

# --- Re-point the active selection, as it ended up after the edits ---
$ws.Range("E16").Select() | Out-Null
